$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1158-1159, shifting the existing data (old rows
# 1158-1196) down by two rows (they become rows 1160-1198).
$ws.Rows("1158:1159").Insert()

# New weekly data row (Primera quality) for fecha 2023-05-29 (serial 45075)
$r = 1158
$ws.Cells.Item($r, 1).Value = 6
$ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($r, 3).Value = "Metropolitana"
$ws.Cells.Item($r, 4).Value = 45075
$ws.Cells.Item($r, 5).Value = 13
$ws.Cells.Item($r, 6).Value = 100112008
$ws.Cells.Item($r, 7).Value = "Coliflor"
$ws.Cells.Item($r, 8).Value = "Sin especificar"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 9200
$ws.Cells.Item($r, 11).Value = 900
$ws.Cells.Item($r, 12).Value = 1000
$ws.Cells.Item($r, 13).Value = 939
$ws.Cells.Item($r, 14).Value = "`$/unidad"
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 939
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = "Hortaliza"

# New weekly data row (Segunda quality) for fecha 2023-05-29 (serial 45075)
$r = 1159
$ws.Cells.Item($r, 1).Value = 6
$ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($r, 3).Value = "Metropolitana"
$ws.Cells.Item($r, 4).Value = 45075
$ws.Cells.Item($r, 5).Value = 13
$ws.Cells.Item($r, 6).Value = 100112008
$ws.Cells.Item($r, 7).Value = "Coliflor"
$ws.Cells.Item($r, 8).Value = "Sin especificar"
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 3300
$ws.Cells.Item($r, 11).Value = 700
$ws.Cells.Item($r, 12).Value = 700
$ws.Cells.Item($r, 13).Value = 700
$ws.Cells.Item($r, 14).Value = "`$/unidad"
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 700
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = "Hortaliza"
